# Update "Forecast Comparison" sheet with a new Week_Start_Date column,
# corrected MyForecast figures, short week labels and boolean holiday flags.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new column B ("Week_Start_Date"); ASIN and the rest shift right.
$ws.Columns.Item(2).Insert()

$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# New week start dates (one per forecast week, rows 2..17)
$weekStartDates = @(
    "2025-01-05",
    "2025-01-12",
    "2025-01-19",
    "2025-01-26",
    "2025-02-02",
    "2025-02-09",
    "2025-02-16",
    "2025-02-23",
    "2025-03-02",
    "2025-03-09",
    "2025-03-16",
    "2025-03-23",
    "2025-03-30",
    "2025-04-06",
    "2025-04-13",
    "2025-04-20"
)

# Corrected MyForecast values (now in column D after the insert)
$myForecast = @(78, 78, 86, 84, 86, 82, 86, 85, 88, 85, 88, 86, 87, 82, 85, 84)

for ($i = 0; $i -lt 16; $i++) {
    $row = $i + 2

    # Week label: "W01" -> "W1" (no leading zero)
    $ws.Cells.Item($row, 1).Value = "W" + ($i + 1)

    # Week_Start_Date column, stored as text (not auto-converted to a date)
    $dateCell = $ws.Cells.Item($row, 2)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $weekStartDates[$i]

    # Corrected MyForecast value (column D)
    $ws.Cells.Item($row, 4).Value = $myForecast[$i]

    # is_holiday_week (column J) becomes a real boolean instead of 0/1 number
    $ws.Cells.Item($row, 10).Value = $false
}

# --- Summary sheet updates -------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)

# All these values are kept as plain text, matching the original cells, so
# force a text number format before assigning to avoid Excel auto-converting
# numeric-looking or date-looking strings into numbers/dates.
$summaryUpdates = @{
    9  = "1352"
    10 = "666"
    11 = "326"
    12 = "88"
    13 = "2025-03-16"
    14 = "78"
    15 = "2025-01-05"
}

foreach ($row in $summaryUpdates.Keys) {
    $cell = $ws2.Cells.Item($row, 2)
    $cell.NumberFormat = "@"
    $cell.Value = $summaryUpdates[$row]
}
